# Apply the "Added sensor, larger opening" edit:
#  - rename two parameter labels (row 20/21): sensorHolderWidth/Height -> detectorOpeningWidth/Height
#    and shrink their values from 3.5 to 3
#  - bump the sensorHolderFlushORingID value (row 26) from 10 to 12
#  - rename fresnelChamfer -> fresnelMirrorChamfer (row 31)
#  - drop the now-orphaned B9 formula (=B8), keeping its cached literal value
#  - move the view/selection down to around row 27

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 31: fresnelChamfer -> fresnelMirrorChamfer (renamed first so the new
#     shared-string entries land in the same order as the authored edit:
#     fresnelMirrorChamfer, detectorOpeningWidth, detectorOpeningHeight) ---
$ws.Range("A31").Value = "fresnelMirrorChamfer"

# --- Row 20: sensorHolderWidth -> detectorOpeningWidth, 3.5 -> 3 ---
$ws.Range("A20").Value = "detectorOpeningWidth"
$ws.Range("B20").Value = 3

# --- Row 21: sensorHolderHeight -> detectorOpeningHeight, 3.5 -> 3 ---
$ws.Range("A21").Value = "detectorOpeningHeight"
$ws.Range("B21").Value = 3

# --- Row 26: sensorHolderFlushORingID value 10 -> 12 ---
$ws.Range("B26").Value = 12

# --- B9 used to be "=B8"; drop the formula but keep the literal cached value ---
$ws.Range("B9").Formula = ""
$ws.Range("B9").Value = 52

# --- Move the selection/view roughly to where row 27 is visible ---
$ws.Range("B27").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 7
$win.ScrollColumn = 1
